$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 42; Excel shifts rows 42:105 down to 43:106
# and extends the used range to A1:R106.
$ws.Rows.Item(42).Insert()

# Populate the newly inserted row 42 with the new weekly price record.
$ws.Range("A42").Value = 8
$ws.Range("B42").Value = "Terminal La Palmera de La Serena"
$ws.Range("C42").Value = "Coquimbo"
$ws.Range("D42").Value = 44546
$ws.Range("E42").Value = 4
$ws.Range("F42").Value = 100112040
$ws.Range("G42").Value = "Cilantro"
$ws.Range("H42").Value = "Sin especificar"
$ws.Range("I42").Value = "Primera"
$ws.Range("J42").Value = 2800
$ws.Range("K42").Value = 2000
$ws.Range("L42").Value = 2500
$ws.Range("M42").Value = 2250
$ws.Range("N42").Value = "`$/atado 1 a 1,5 kilos"
$ws.Range("O42").Value = "Provincia del Elquí"
$ws.Range("P42").Value = 1500
$ws.Range("Q42").Value = 1.5
$ws.Range("R42").Value = "Hortaliza"
